$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H7").Value = 500
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 500
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 500
$ws.Range("M7").ClearContents()
$ws.Range("N7").Value = -724

$ws.Range("H14").Value = 500
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 500
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 500
$ws.Range("M14").ClearContents()
$ws.Range("N14").Value = -882

$ws.Range("H17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("N17").ClearContents()

$ws.Range("H74").Value = 2186.125
$ws.Range("I74").Value = 1784.1428
$ws.Range("K74").Value = 1784.1428
$ws.Range("M74").Value = -848.1428000000001

$ws.Range("H77").Value = 2186.125
$ws.Range("I77").Value = 1784.1428
$ws.Range("K77").Value = 8920.714
$ws.Range("M77").Value = -4240.714

$ws.Range("H132").Value = 1808.5
$ws.Range("I132").Value = 2208.1667
$ws.Range("J132").Value = 609.5
$ws.Range("K132").Value = 6624.500100000001
$ws.Range("L132").Value = 1828.5
$ws.Range("M132").Value = -4094.500100000001
$ws.Range("N132").Value = -6888.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 1003.3333
$ws.Range("I110").Value = 1003.3333
$ws.Range("K110").Value = 1003.3333
$ws.Range("M110").Value = 1041.6667

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1338
$ws.Range("I86").Value = 1000
$ws.Range("K86").Value = 1000
$ws.Range("M86").Value = 123

$ws.Range("H89").Value = 1338
$ws.Range("I89").Value = 1000
$ws.Range("K89").Value = 5000
$ws.Range("M89").Value = 616

$ws.Range("H134").Value = 5081.8887
$ws.Range("I134").Value = 1787.1666
$ws.Range("J134").Value = 11671.333
$ws.Range("K134").Value = 5361.4998
$ws.Range("L134").Value = 35013.999
$ws.Range("M134").Value = -2826.4998
$ws.Range("N134").Value = -40083.999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 198
$ws.Range("I22").Value = 196.8
$ws.Range("J22").Value = 200
$ws.Range("K22").Value = 196.8
$ws.Range("L22").Value = 200
$ws.Range("M22").Value = 153.2
$ws.Range("N22").Value = -900

$ws.Range("H35").Value = 2017
$ws.Range("I35").Value = 2017
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 2017
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = -1723
$ws.Range("N35").ClearContents()

$ws.Range("H69").Value = 27599
$ws.Range("I69").Value = 13995
$ws.Range("K69").Value = 13995
$ws.Range("M69").Value = -13246

$ws.Range("H72").Value = 27599
$ws.Range("I72").Value = 13995
$ws.Range("K72").Value = 41985
$ws.Range("M72").Value = -38241

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H35").Value = 1002.5
$ws.Range("I35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("M35").ClearContents()

$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("M64").ClearContents()

$ws.Range("H67").Value = 0
$ws.Range("I67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("M67").ClearContents()

$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()

$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()

$ws.Range("H111").Value = 1000
$ws.Range("I111").Value = 1000
$ws.Range("K111").Value = 3000
$ws.Range("M111").Value = 67

$ws.Range("H116").Value = 883
$ws.Range("I116").Value = 883
$ws.Range("K116").Value = 2649
$ws.Range("M116").Value = 793

$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()

$ws.Range("I120").Value = 10000
$ws.Range("J120").Value = 0
$ws.Range("K120").Value = 30000
$ws.Range("L120").Value = 0
$ws.Range("M120").Value = -25162
$ws.Range("N120").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 85.875
$ws.Range("I2").Value = 29
$ws.Range("K2").Value = 29
$ws.Range("M2").Value = 84

$ws.Range("H126").Value = 2162
$ws.Range("I126").Value = 2162
$ws.Range("K126").Value = 6486
$ws.Range("M126").Value = -4016

$ws.Range("H132").Value = 5298.75
$ws.Range("I132").Value = 4358.5
$ws.Range("J132").Value = 10000
$ws.Range("K132").Value = 13075.5
$ws.Range("L132").Value = 30000
$ws.Range("M132").Value = -10545.5
$ws.Range("N132").Value = -35060

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("N13").ClearContents()

$ws.Range("H46").Value = 7450
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()

$ws.Range("H58").Value = 5000
$ws.Range("I58").Value = 5000
$ws.Range("K58").Value = 5000
$ws.Range("M58").Value = -4740

$ws.Range("H62").Value = 25000
$ws.Range("J62").Value = 25000
$ws.Range("L62").Value = 25000
$ws.Range("N62").Value = -26248

$ws.Range("H65").Value = 25000
$ws.Range("J65").Value = 25000
$ws.Range("L65").Value = 75000
$ws.Range("N65").Value = -81240

$ws.Range("H93").Value = 2566
$ws.Range("I93").Value = 2566
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 2566
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = -1318
$ws.Range("N93").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H39").Value = 0
$ws.Range("I39").Value = 0
$ws.Range("K39").Value = 0
$ws.Range("M39").ClearContents()

$ws.Range("H51").Value = 60000
$ws.Range("I51").Value = 60000
$ws.Range("K51").Value = 60000
$ws.Range("M51").Value = -59490

$ws.Range("H52").Value = 18995
$ws.Range("I52").Value = 18995
$ws.Range("K52").Value = 18995
$ws.Range("M52").Value = -18769

$ws.Range("H54").Value = 31449.5
$ws.Range("J54").Value = 31449.5
$ws.Range("L54").Value = 31449.5
$ws.Range("N54").Value = -32489.5

$ws.Range("H55").Value = 0
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("M55").ClearContents()
$ws.Range("N55").ClearContents()

$ws.Range("H81").Value = 199.5
$ws.Range("I81").Value = 199.5
$ws.Range("K81").Value = 399
$ws.Range("M81").Value = 662

$ws.Range("H84").Value = 199.5
$ws.Range("I84").Value = 199.5
$ws.Range("K84").Value = 1995
$ws.Range("M84").Value = 3309
